$wb = $excel.ActiveWorkbook

# Sheet: ALC
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1172.0769
$ws.Range("I19").Value = 1002.5
$ws.Range("J19").Value = 1247.4445
$ws.Range("K19").Value = 1002.5
$ws.Range("L19").Value = 1247.4445
$ws.Range("M19").Value = -827.5
$ws.Range("N19").Value = -1597.4445
$ws.Range("H75").Value = 28750
$ws.Range("J75").Value = 28750
$ws.Range("L75").Value = 28750
$ws.Range("N75").Value = -30622
$ws.Range("H78").Value = 28750
$ws.Range("J78").Value = 28750
$ws.Range("L78").Value = 86250
$ws.Range("N78").Value = -95610
$ws.Range("H113").Value = 78720.46
$ws.Range("J113").Value = 1856.7778
$ws.Range("L113").Value = 1856.7778
$ws.Range("N113").Value = -8364.7778
$ws.Range("H130").Value = 30380
$ws.Range("J130").Value = 30380
$ws.Range("L130").Value = 30380
$ws.Range("N130").Value = -40420
$ws.Range("H138").Value = 5989.9014
$ws.Range("I138").Value = 1871.1578
$ws.Range("J138").Value = 7494.827
$ws.Range("K138").Value = 5613.4734
$ws.Range("L138").Value = 22484.481
$ws.Range("M138").Value = -473.4733999999999
$ws.Range("N138").Value = -32764.481

# Sheet: ARM
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()
$ws.Range("H32").Value = 22459.902
$ws.Range("I32").Value = 4081.8381
$ws.Range("K32").Value = 4081.8381
$ws.Range("M32").Value = -3794.8381

# Sheet: BSM
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 103
$ws.Range("I7").Value = 103
$ws.Range("K7").Value = 103
$ws.Range("M7").Value = 10
$ws.Range("H19").Value = 22500
$ws.Range("J19").Value = 22500
$ws.Range("L19").Value = 22500
$ws.Range("N19").Value = -22846
$ws.Range("H99").Value = 2040
$ws.Range("I99").Value = 1786.6666
$ws.Range("J99").Value = 2800
$ws.Range("K99").Value = 1786.6666
$ws.Range("L99").Value = 2800
$ws.Range("M99").Value = -288.6666
$ws.Range("N99").Value = -5796

# Sheet: CRP
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2619.2
$ws.Range("J62").Value = 2619.2
$ws.Range("L62").Value = 2619.2
$ws.Range("N62").Value = -3867.2
$ws.Range("H65").Value = 2619.2
$ws.Range("J65").Value = 2619.2
$ws.Range("L65").Value = 13096
$ws.Range("N65").Value = -19336
$ws.Range("H76").Value = 505020
$ws.Range("I76").Value = 505020
$ws.Range("K76").Value = 505020
$ws.Range("M76").Value = -504705
$ws.Range("H79").Value = 505020
$ws.Range("I79").Value = 505020
$ws.Range("K79").Value = 505020
$ws.Range("M79").Value = -503928
$ws.Range("H99").Value = 23314.5
$ws.Range("I99").Value = 6964.5
$ws.Range("J99").Value = 56014.5
$ws.Range("K99").Value = 6964.5
$ws.Range("L99").Value = 56014.5
$ws.Range("M99").Value = -5466.5
$ws.Range("N99").Value = -59010.5
$ws.Range("H105").Value = 2639
$ws.Range("I105").Value = 2772.375
$ws.Range("K105").Value = 2772.375
$ws.Range("M105").Value = -1025.375
$ws.Range("H126").Value = 23314.5
$ws.Range("I126").Value = 6964.5
$ws.Range("J126").Value = 56014.5
$ws.Range("K126").Value = 20893.5
$ws.Range("L126").Value = 168043.5
$ws.Range("M126").Value = -18423.5
$ws.Range("N126").Value = -172983.5

# Sheet: CUL
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 19.5625
$ws.Range("I12").Value = 11.333333
$ws.Range("J12").Value = 24.5
$ws.Range("K12").Value = 33.999999
$ws.Range("L12").Value = 73.5
$ws.Range("M12").Value = 139.000001
$ws.Range("N12").Value = -419.5
$ws.Range("H22").Value = 18160
$ws.Range("J22").Value = 26990
$ws.Range("L22").Value = 80970
$ws.Range("N22").Value = -81308
$ws.Range("H23").Value = 1179.9412
$ws.Range("I23").Value = 2246
$ws.Range("J23").Value = 735.75
$ws.Range("K23").Value = 6738
$ws.Range("L23").Value = 2207.25
$ws.Range("M23").Value = -6503
$ws.Range("N23").Value = -2677.25
$ws.Range("H27").Value = 18160
$ws.Range("J27").Value = 26990
$ws.Range("L27").Value = 80970
$ws.Range("N27").Value = -81174
$ws.Range("H47").Value = 135.9
$ws.Range("I47").Value = 119.875
$ws.Range("K47").Value = 359.625
$ws.Range("M47").Value = 71.375
$ws.Range("H50").Value = 2000.5555
$ws.Range("I50").Value = 2747.25
$ws.Range("J50").Value = 1403.2
$ws.Range("K50").Value = 8241.75
$ws.Range("L50").Value = 4209.6
$ws.Range("M50").Value = -7760.75
$ws.Range("N50").Value = -5171.6
$ws.Range("H53").Value = 2000.5555
$ws.Range("I53").Value = 2747.25
$ws.Range("J53").Value = 1403.2
$ws.Range("K53").Value = 8241.75
$ws.Range("L53").Value = 4209.6
$ws.Range("M53").Value = -7760.75
$ws.Range("N53").Value = -5171.6
$ws.Range("H70").Value = 145073.14
$ws.Range("I70").Value = 501256
$ws.Range("J70").Value = 2600
$ws.Range("K70").Value = 1503768
$ws.Range("L70").Value = 7800
$ws.Range("M70").Value = -1503453
$ws.Range("N70").Value = -8430
$ws.Range("H73").Value = 145073.14
$ws.Range("I73").Value = 501256
$ws.Range("J73").Value = 2600
$ws.Range("K73").Value = 1503768
$ws.Range("L73").Value = 7800
$ws.Range("M73").Value = -1502676
$ws.Range("N73").Value = -9984
$ws.Range("H131").Value = 826.09186
$ws.Range("J131").Value = 869.1705
$ws.Range("L131").Value = 2607.5115
$ws.Range("N131").Value = -12687.5115

# Sheet: GSM
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 55002.75
$ws.Range("H70").Value = 254502.12
$ws.Range("I70").Value = 503002
$ws.Range("J70").Value = 6002.25
$ws.Range("K70").Value = 503002
$ws.Range("L70").Value = 6002.25
$ws.Range("M70").Value = -502732
$ws.Range("N70").Value = -6542.25
$ws.Range("H73").Value = 254502.12
$ws.Range("I73").Value = 503002
$ws.Range("J73").Value = 6002.25
$ws.Range("K73").Value = 503002
$ws.Range("L73").Value = 6002.25
$ws.Range("M73").Value = -502066
$ws.Range("N73").Value = -7874.25
$ws.Range("H102").Value = 2994.9285
$ws.Range("J102").Value = 3555.6
$ws.Range("L102").Value = 3555.6
$ws.Range("N102").Value = -6799.6

# Sheet: LTW
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2720
$ws.Range("I61").Value = 900
$ws.Range("J61").Value = 3326.6667
$ws.Range("K61").Value = 900
$ws.Range("L61").Value = 3326.6667
$ws.Range("M61").Value = -698
$ws.Range("N61").Value = -3730.6667
$ws.Range("H113").Value = 2720
$ws.Range("I113").Value = 900
$ws.Range("J113").Value = 3326.6667
$ws.Range("K113").Value = 900
$ws.Range("L113").Value = 3326.6667
$ws.Range("M113").Value = 1270
$ws.Range("N113").Value = -7666.6667
$ws.Range("H132").Value = 3607.4614
$ws.Range("I132").Value = 3671.76
$ws.Range("J132").Value = 2000
$ws.Range("K132").Value = 11015.28
$ws.Range("L132").Value = 6000
$ws.Range("M132").Value = -8485.28
$ws.Range("N132").Value = -11060

# Sheet: WVR
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H29").Value = 70011
$ws.Range("I29").Value = 0
$ws.Range("K29").Value = 0
$ws.Range("M29").ClearContents()
$ws.Range("H62").Value = 5685845.5
$ws.Range("I62").Value = 31251200
$ws.Range("J62").Value = 4655.5557
$ws.Range("K62").Value = 31251200
$ws.Range("L62").Value = 4655.5557
$ws.Range("M62").Value = -31250576
$ws.Range("N62").Value = -5903.5557
$ws.Range("H65").Value = 5685845.5
$ws.Range("I65").Value = 31251200
$ws.Range("J65").Value = 4655.5557
$ws.Range("K65").Value = 156256000
$ws.Range("L65").Value = 23277.7785
$ws.Range("M65").Value = -156252880
$ws.Range("N65").Value = -29517.7785
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("M122").ClearContents()
$ws.Range("N122").ClearContents()

